$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.654.87'
$ws.Range("E2").Value = '  -0.29%  '
$ws.Range("D3").Value = '2.111.88'
$ws.Range("E3").Value = '  +9.56%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '253.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.664'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.51%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '48.30'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +9.23%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '59.50'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.15%  '
$ws.Range("E10").Value = '  +0.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0746'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.92%  '
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("D13").Value = '2.421.12'
$ws.Range("E13").Value = '  +9.66%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.31'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.832'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.02%  '
$ws.Range("D16").Value = '2.109.80'
$ws.Range("E16").Value = '  +9.56%  '
$ws.Range("E17").Value = '  -0.61%  '
$ws.Range("D18").Value = '36.646.10'
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.50'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.62%  '
$ws.Range("E20").Value = '  -3.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '240.54'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.71%  '
$ws.Range("E23").Value = '  -1.29%  '
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("E25").Value = '  -7.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '172.72'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.57'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +14.54%  '
$ws.Range("E28").Value = '  +3.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.02'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -9.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '29.74'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +66.90%  '
$ws.Range("E31").Value = '  -4.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.49'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.97%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0604'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0914'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.964'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.65%  '
$ws.Range("E36").Value = '  +14.00%  '
$ws.Range("B37").Value = 'BinanceUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("B38").Value = 'WEMIXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.89'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.31%  '
$ws.Range("E39").Value = '  -6.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.34'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -12.05%  '
$ws.Range("E41").Value = '  +6.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0226'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '98.66'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.77'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.03'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.61%  '
$ws.Range("D46").Value = '1.345.70'
$ws.Range("E46").Value = '  +0.09%  '
$ws.Range("E47").Value = '  +11.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0843'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.24%  '
$ws.Range("D49").Value = '2.304.73'
$ws.Range("E49").Value = '  +9.19%  '
$ws.Range("E50").Value = '  +1.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.28'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.15%  '
